# Capitalize the first letter of the phenotype names in column B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Endometriosis"
$ws.Range("B3").Value = "Polycystic ovary syndrome"
$ws.Range("B4").Value = "Polycystic ovary syndrome"
$ws.Range("B5").Value = "Recurrent spontaneous abortion"
$ws.Range("B6").Value = "Recurrent spontaneous abortion"
$ws.Range("B7").Value = "Recurrent spontaneous abortion"
$ws.Range("B12").Value = "Oligoasthenoteratozoospermia"
$ws.Range("B13").Value = "Oligoasthenoteratozoospermia"
$ws.Range("B14").Value = "Non-obstructive azoospermia"
